$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header labels in row 1 (shared strings) ---
$ws.Range("C1").Value = "AN_ANUAL_202105"
$ws.Range("D1").Value = "AN_SEM_202105"
$ws.Range("E1").Value = "AN_TRI_202105"
$ws.Range("F1").Value = "AC_202105"
$ws.Range("G1").Value = "AN_202105"
$ws.Range("H1").Value = "SPI1_202105"
$ws.Range("I1").Value = "SPI3_202105"
$ws.Range("J1").Value = "SPI6_202105"
$ws.Range("K1").Value = "SPI12_202105"
$ws.Range("L1").Value = "AN_COTA_202103"
$ws.Range("M1").Value = "AN_COTA_202104"
$ws.Range("N1").Value = "AN_COTA_20210531"

# --- Update numeric data cells rows 2-42 ---
# Row 2
$ws.Range("C2").Value = -22.4
$ws.Range("D2").Value = -17.4
$ws.Range("E2").Value = -19.2
$ws.Range("F2").Value = 115.1
$ws.Range("G2").Value = -5.1
$ws.Range("H2").Value = 0.09
$ws.Range("I2").Value = -0.37
$ws.Range("J2").Value = -0.36
$ws.Range("K2").Value = -0.63
$ws.Range("L2").Value = -1.85
$ws.Range("M2").Value = -3.18
$ws.Range("N2").Value = -7.16
# Row 3
$ws.Range("C3").Value = -9.1
$ws.Range("D3").Value = -0.4
$ws.Range("E3").Value = -17.5
$ws.Range("F3").Value = 92.3
$ws.Range("G3").Value = 17.3
$ws.Range("H3").Value = -0.05
$ws.Range("I3").Value = -0.8
$ws.Range("J3").Value = -0.19
$ws.Range("K3").Value = -0.38
$ws.Range("L3").Value = 18.97
$ws.Range("M3").Value = -31.61
$ws.Range("N3").Value = -12.35
# Row 4
$ws.Range("C4").Value = -30.1
$ws.Range("D4").Value = -23.7
$ws.Range("E4").Value = -43.3
$ws.Range("F4").Value = 107
$ws.Range("G4").Value = -19
$ws.Range("H4").Value = -0.22
$ws.Range("I4").Value = -0.51
$ws.Range("J4").Value = -0.09
$ws.Range("K4").Value = -0.47
# Row 5
$ws.Range("C5").Value = -7.1
$ws.Range("D5").Value = 12.9
$ws.Range("E5").Value = -62.1
$ws.Range("F5").Value = 33.8
$ws.Range("G5").Value = -74.2
$ws.Range("H5").Value = -1.28
$ws.Range("I5").Value = -3.47
$ws.Range("J5").Value = -0.91
$ws.Range("K5").Value = -1.1
$ws.Range("L5").Value = -37.18
$ws.Range("M5").Value = -47.57
$ws.Range("N5").Value = -63
# Row 6
$ws.Range("C6").Value = -25.6
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = -63.7
$ws.Range("F6").Value = 24.3
$ws.Range("G6").Value = -80.8
$ws.Range("H6").Value = -1.58
$ws.Range("I6").Value = -2.35
$ws.Range("J6").Value = -0.93
$ws.Range("K6").Value = -1.86
$ws.Range("L6").Value = 9.01
$ws.Range("M6").Value = -3.55
$ws.Range("N6").Value = -11.81
# Row 7
$ws.Range("C7").Value = -4.8
$ws.Range("D7").Value = 18.2
$ws.Range("E7").Value = -45.7
$ws.Range("F7").Value = 35.8
$ws.Range("G7").Value = -75.2
$ws.Range("H7").Value = -1.29
$ws.Range("I7").Value = -1.69
$ws.Range("J7").Value = -0.73
$ws.Range("K7").Value = -1.33
# Row 8
$ws.Range("C8").Value = -23.3
$ws.Range("D8").Value = -20.1
$ws.Range("E8").Value = -31.5
$ws.Range("F8").Value = 80
$ws.Range("G8").Value = -23.8
$ws.Range("H8").Value = 0.22
$ws.Range("I8").Value = -0.4
$ws.Range("J8").Value = 0.21
$ws.Range("K8").Value = 0.01
$ws.Range("L8").Value = 10.54
$ws.Range("M8").Value = -3.26
$ws.Range("N8").Value = -2.48
# Row 9
$ws.Range("C9").Value = -30.5
$ws.Range("D9").Value = -22.2
$ws.Range("E9").Value = -73.4
$ws.Range("F9").Value = 24.4
$ws.Range("G9").Value = -84
$ws.Range("H9").Value = -1.27
$ws.Range("I9").Value = -2.11
$ws.Range("J9").Value = -0.56
$ws.Range("K9").Value = -0.98
# Row 10
$ws.Range("C10").Value = -24.1
$ws.Range("D10").Value = -17.6
$ws.Range("E10").Value = -28.9
$ws.Range("F10").Value = 100.1
$ws.Range("G10").Value = 15.5
$ws.Range("H10").Value = 0.58
$ws.Range("I10").Value = -0.52
$ws.Range("J10").Value = -0.46
$ws.Range("K10").Value = -0.5
$ws.Range("L10").Value = -18.36
$ws.Range("M10").Value = -21.81
$ws.Range("N10").Value = 9.38
# Row 11
$ws.Range("C11").Value = -18.5
$ws.Range("D11").Value = -11.2
$ws.Range("E11").Value = -8.8
$ws.Range("F11").Value = 104.6
$ws.Range("G11").Value = 1.3
$ws.Range("H11").Value = 0.34
$ws.Range("I11").Value = -0.17
$ws.Range("J11").Value = -0.12
$ws.Range("K11").Value = -0.23
$ws.Range("L11").Value = 37.37
$ws.Range("M11").Value = -4.65
$ws.Range("N11").Value = -6.15
# Row 12
$ws.Range("C12").Value = -18.8
$ws.Range("D12").Value = -12
$ws.Range("E12").Value = -10.5
$ws.Range("F12").Value = 108.2
$ws.Range("G12").Value = -0.8
$ws.Range("H12").Value = 0.39
$ws.Range("I12").Value = -0.08
$ws.Range("J12").Value = -0.15
$ws.Range("K12").Value = -0.2
$ws.Range("L12").Value = 37.37
$ws.Range("M12").Value = -4.65
$ws.Range("N12").Value = -6.15
# Row 13
$ws.Range("C13").Value = -38.9
$ws.Range("D13").Value = -38.4
$ws.Range("E13").Value = -79.5
$ws.Range("F13").Value = 19.1
$ws.Range("G13").Value = -88.7
$ws.Range("H13").Value = -1.26
$ws.Range("I13").Value = -2.55
$ws.Range("J13").Value = -0.57
$ws.Range("K13").Value = -0.98
# Row 14
$ws.Range("C14").Value = -27
$ws.Range("D14").Value = -21
$ws.Range("E14").Value = -23.2
$ws.Range("F14").Value = 91.5
$ws.Range("G14").Value = -19.1
$ws.Range("H14").Value = 0.24
$ws.Range("I14").Value = -0.13
$ws.Range("J14").Value = 0.29
$ws.Range("K14").Value = -0.03
$ws.Range("L14").Value = 20.35
$ws.Range("M14").Value = -20.59
$ws.Range("N14").Value = -18.42
# Row 15
$ws.Range("C15").Value = -23.6
$ws.Range("D15").Value = -14.9
$ws.Range("E15").Value = -51.3
$ws.Range("F15").Value = 71.6
$ws.Range("G15").Value = -51.6
$ws.Range("H15").Value = -1.03
$ws.Range("I15").Value = -2.06
$ws.Range("J15").Value = -0.73
$ws.Range("K15").Value = -1.06
$ws.Range("L15").ClearContents()
$ws.Range("M15").ClearContents()
$ws.Range("N15").ClearContents()
# Row 16
$ws.Range("C16").Value = -24.7
$ws.Range("D16").Value = -16.9
$ws.Range("E16").Value = -51.6
$ws.Range("F16").Value = 70.9
$ws.Range("G16").Value = -60.1
$ws.Range("H16").Value = -1.25
$ws.Range("I16").Value = -2.56
$ws.Range("J16").Value = -0.61
$ws.Range("K16").Value = -1.01
$ws.Range("L16").Value = -15.97
$ws.Range("M16").Value = -13.16
$ws.Range("N16").Value = -28.32
# Row 17
$ws.Range("C17").Value = -45.8
$ws.Range("D17").Value = -42.6
$ws.Range("E17").Value = -47.8
$ws.Range("F17").Value = 65.1
$ws.Range("G17").Value = -53.1
$ws.Range("H17").Value = -0.26
$ws.Range("I17").Value = -0.61
$ws.Range("J17").Value = -1.3
$ws.Range("K17").Value = -0.92
# Row 18
$ws.Range("C18").Value = -37.6
$ws.Range("D18").Value = -34.7
$ws.Range("E18").Value = -32
$ws.Range("F18").Value = 117.6
$ws.Range("G18").Value = -22.5
$ws.Range("H18").Value = 0.09
$ws.Range("I18").Value = -0.37
$ws.Range("J18").Value = -0.36
$ws.Range("K18").Value = -0.63
$ws.Range("L18").Value = -1.85
$ws.Range("M18").Value = -3.18
$ws.Range("N18").Value = -7.16
# Row 19
$ws.Range("C19").Value = -8.8
$ws.Range("D19").Value = 5.4
$ws.Range("E19").Value = -50.5
$ws.Range("F19").Value = 48.5
$ws.Range("G19").Value = -62
$ws.Range("H19").Value = -1.33
$ws.Range("I19").Value = -2.14
$ws.Range("J19").Value = -0.78
$ws.Range("K19").Value = -1.27
# Row 20
$ws.Range("C20").Value = -33.3
$ws.Range("D20").Value = -28.6
$ws.Range("E20").Value = -47.7
$ws.Range("F20").Value = 95.6
$ws.Range("G20").Value = -38.4
$ws.Range("H20").Value = -1.16
$ws.Range("I20").Value = -1.07
$ws.Range("J20").Value = 0.48
$ws.Range("K20").Value = -0.05
# Row 21
$ws.Range("C21").Value = -24.6
$ws.Range("D21").Value = -18.2
$ws.Range("E21").Value = -50.7
$ws.Range("F21").Value = 73.8
$ws.Range("G21").Value = -57.4
$ws.Range("H21").Value = -1.23
$ws.Range("I21").Value = -2.49
$ws.Range("J21").Value = -0.61
$ws.Range("K21").Value = -1.01
$ws.Range("L21").Value = 752.44
$ws.Range("M21").Value = 706.31
$ws.Range("N21").Value = -8.7
# Row 22
$ws.Range("C22").Value = -3.9
$ws.Range("D22").Value = 12.8
$ws.Range("E22").Value = -35.3
$ws.Range("F22").Value = 58.4
$ws.Range("G22").Value = -64.6
$ws.Range("H22").Value = -1.15
$ws.Range("I22").Value = -1.82
$ws.Range("J22").Value = -0.7
$ws.Range("K22").Value = -1.03
# Row 23
$ws.Range("C23").Value = -21.7
$ws.Range("D23").Value = -0.8
$ws.Range("E23").Value = -64.1
$ws.Range("F23").Value = 30.7
$ws.Range("G23").Value = -80.7
$ws.Range("H23").Value = -1.37
$ws.Range("I23").Value = -2.77
$ws.Range("J23").Value = -2.11
$ws.Range("K23").Value = -1.9
# Row 24
$ws.Range("C24").Value = -31.3
$ws.Range("D24").Value = -18.3
$ws.Range("E24").Value = -71.9
$ws.Range("F24").Value = 24.4
$ws.Range("G24").Value = -84.9
$ws.Range("H24").Value = -1.26
$ws.Range("I24").Value = -2.7
$ws.Range("J24").Value = -0.47
$ws.Range("K24").Value = -0.96
# Row 25
$ws.Range("C25").Value = -22.9
$ws.Range("D25").Value = -9.9
$ws.Range("E25").Value = -27
$ws.Range("F25").Value = 71.9
$ws.Range("G25").Value = -36.6
$ws.Range("H25").Value = -0.13
$ws.Range("I25").Value = -1.11
$ws.Range("J25").Value = -0.2
$ws.Range("K25").Value = -0.41
# Row 26
$ws.Range("C26").Value = -29.3
$ws.Range("D26").Value = -20.1
$ws.Range("E26").Value = -33.8
$ws.Range("F26").Value = 77.6
$ws.Range("G26").Value = -27.5
$ws.Range("H26").Value = -0.02
$ws.Range("I26").Value = -0.99
$ws.Range("J26").Value = -0.48
$ws.Range("K26").Value = -0.87
# Row 27
$ws.Range("C27").Value = -46.1
$ws.Range("D27").Value = -44.8
$ws.Range("E27").Value = -48.9
$ws.Range("F27").Value = 74.9
$ws.Range("G27").Value = -44.9
$ws.Range("H27").Value = -0.13
$ws.Range("I27").Value = -0.26
$ws.Range("J27").Value = -0.61
$ws.Range("K27").Value = -0.67
# Row 28
$ws.Range("C28").Value = -9.2
$ws.Range("D28").Value = 11
$ws.Range("E28").Value = -46.1
$ws.Range("F28").Value = 47.3
$ws.Range("G28").Value = -71.7
$ws.Range("H28").Value = -1.27
$ws.Range("I28").Value = -2.9
$ws.Range("J28").Value = -0.63
$ws.Range("K28").Value = -1.02
$ws.Range("L28").Value = -37.18
$ws.Range("M28").Value = -47.57
$ws.Range("N28").Value = -63
# Row 29
$ws.Range("C29").Value = -21.7
$ws.Range("D29").Value = -7.7
$ws.Range("E29").Value = -18.9
$ws.Range("F29").Value = 97.6
$ws.Range("G29").Value = -11.1
$ws.Range("H29").Value = 0.4
$ws.Range("I29").Value = 0.01
$ws.Range("J29").Value = 0.32
$ws.Range("K29").Value = 0.19
$ws.Range("L29").Value = 84.34
$ws.Range("M29").Value = -12.83
$ws.Range("N29").Value = -24.54
# Row 30
$ws.Range("C30").Value = -11.3
$ws.Range("D30").Value = 3.6
$ws.Range("E30").Value = -37.4
$ws.Range("F30").Value = 74.6
$ws.Range("G30").Value = -55.4
$ws.Range("H30").Value = -1.25
$ws.Range("I30").Value = -2.62
$ws.Range("J30").Value = -0.64
$ws.Range("K30").Value = -1.03
$ws.Range("L30").Value = -37.18
$ws.Range("M30").Value = -47.57
$ws.Range("N30").Value = -63
# Row 31
$ws.Range("C31").Value = -18.8
$ws.Range("D31").Value = -6.1
$ws.Range("E31").Value = -50.7
$ws.Range("F31").Value = 51.4
$ws.Range("G31").Value = -70.2
$ws.Range("H31").Value = -1.26
$ws.Range("I31").Value = -2.61
$ws.Range("J31").Value = -0.57
$ws.Range("K31").Value = -0.99
# Row 32
$ws.Range("C32").Value = -9.9
$ws.Range("D32").Value = 10
$ws.Range("E32").Value = -46.1
$ws.Range("F32").Value = 47.2
$ws.Range("G32").Value = -71.7
$ws.Range("H32").Value = -1.27
$ws.Range("I32").Value = -2.85
$ws.Range("J32").Value = -0.61
$ws.Range("K32").Value = -1.01
$ws.Range("L32").Value = -37.18
$ws.Range("M32").Value = -47.57
$ws.Range("N32").Value = -63
# Row 33
$ws.Range("C33").Value = -26.5
$ws.Range("D33").Value = -7.3
$ws.Range("E33").Value = -65.3
$ws.Range("F33").Value = 26.7
$ws.Range("G33").Value = -82.5
$ws.Range("H33").Value = -1.36
$ws.Range("I33").Value = -2.7
$ws.Range("J33").Value = -2.41
$ws.Range("K33").Value = -2.04
# Row 34
$ws.Range("C34").Value = -25.4
$ws.Range("D34").Value = -19.7
$ws.Range("E34").Value = -26.7
$ws.Range("F34").Value = 90.1
$ws.Range("G34").Value = -18.2
$ws.Range("H34").Value = 0.21
$ws.Range("I34").Value = -0.14
$ws.Range("J34").Value = 0.31
$ws.Range("K34").Value = -0.05
$ws.Range("L34").Value = 20.35
$ws.Range("M34").Value = -20.59
$ws.Range("N34").Value = -18.42
# Row 35
$ws.Range("C35").Value = -50.7
$ws.Range("D35").Value = -51.4
$ws.Range("E35").Value = -55
$ws.Range("F35").Value = 70.7
$ws.Range("G35").Value = -49.7
$ws.Range("H35").Value = -0.3
$ws.Range("I35").Value = -0.28
$ws.Range("J35").Value = -0.96
$ws.Range("K35").Value = -0.79
# Row 36
$ws.Range("C36").Value = -31.8
$ws.Range("D36").Value = -25.8
$ws.Range("E36").Value = -25.9
$ws.Range("F36").Value = 95.9
$ws.Range("G36").Value = -33.9
$ws.Range("H36").Value = 0.26
$ws.Range("I36").Value = -0.02
$ws.Range("J36").Value = 0.69
$ws.Range("K36").Value = 0.33
$ws.Range("L36").Value = 15.44
$ws.Range("M36").Value = -58.57
$ws.Range("N36").Value = -61.97
# Row 37
$ws.Range("C37").Value = -9.1
$ws.Range("D37").Value = 5.3
$ws.Range("E37").Value = -34
$ws.Range("F37").Value = 62.6
$ws.Range("G37").Value = -64
$ws.Range("H37").Value = -1.27
$ws.Range("I37").Value = -2.84
$ws.Range("J37").Value = -0.69
$ws.Range("K37").Value = -1.04
$ws.Range("L37").Value = -37.18
$ws.Range("M37").Value = -47.57
$ws.Range("N37").Value = -63
# Row 38
$ws.Range("C38").Value = -20.6
$ws.Range("D38").Value = -16
$ws.Range("E38").Value = -33.5
$ws.Range("F38").Value = 91.2
$ws.Range("G38").Value = -1.8
$ws.Range("H38").Value = 0.57
$ws.Range("I38").Value = -0.49
$ws.Range("J38").Value = -0.44
$ws.Range("K38").Value = -0.47
$ws.Range("L38").Value = -18.36
$ws.Range("M38").Value = -21.81
$ws.Range("N38").Value = 9.38
# Row 39
$ws.Range("C39").Value = -8.3
$ws.Range("D39").Value = 7.4
$ws.Range("E39").Value = -32.4
$ws.Range("F39").Value = 77.1
$ws.Range("G39").Value = -53.4
$ws.Range("H39").Value = -1.25
$ws.Range("I39").Value = -2.59
$ws.Range("J39").Value = -0.65
$ws.Range("K39").Value = -1.03
$ws.Range("L39").Value = -37.18
$ws.Range("M39").Value = -47.57
$ws.Range("N39").Value = -63
# Row 40
$ws.Range("C40").Value = -8.6
$ws.Range("D40").Value = 7.3
$ws.Range("E40").Value = -38.8
$ws.Range("F40").Value = 59.8
$ws.Range("G40").Value = -64.8
$ws.Range("H40").Value = -1.27
$ws.Range("I40").Value = -2.81
$ws.Range("J40").Value = -0.65
$ws.Range("K40").Value = -1.02
$ws.Range("L40").Value = -37.18
$ws.Range("M40").Value = -47.57
$ws.Range("N40").Value = -63
# Row 41
$ws.Range("C41").Value = -21.7
$ws.Range("D41").Value = -14.5
$ws.Range("E41").Value = -18.7
$ws.Range("F41").Value = 92.9
$ws.Range("G41").Value = -14.7
$ws.Range("H41").Value = 0.03
$ws.Range("I41").Value = -0.42
$ws.Range("J41").Value = 0.15
$ws.Range("K41").Value = -0.23
$ws.Range("L41").Value = 20.35
$ws.Range("M41").Value = -20.59
$ws.Range("N41").Value = -18.42
# Row 42
$ws.Range("C42").Value = -26.2
$ws.Range("D42").Value = -1.4
$ws.Range("E42").Value = -64.9
$ws.Range("F42").Value = 26.9
$ws.Range("G42").Value = -81.2
$ws.Range("H42").Value = -1.58
$ws.Range("I42").Value = -2.36
$ws.Range("J42").Value = -0.94
$ws.Range("K42").Value = -1.87
$ws.Range("L42").Value = 9.01
$ws.Range("M42").Value = -3.55
$ws.Range("N42").Value = -11.81
